# Add a new user row ("aaa" in the Username column) to the Users sheet,
# matching the upstream commit that appended a row under the existing
# Username/Email/Password/Avatar header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "aaa"

# Leave the cursor where the author's saved file shows it (cell G7),
# which is recorded in the sheet's stored selection.
$ws.Range("G7").Select() | Out-Null
